# "1st changes of mifos to finflux"
#
# The "Repayment schedule" sheet becomes the active sheet (it was "Transactions"
# before), a new blank column is inserted before column N (shifting the old
# N/O/P columns - "Late", "heading"/blank, "Outstanding" - one place to the
# right, to O/P/Q), and the active cell selection on that sheet moves to K13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active/selected sheet (was "Transactions").
$ws.Activate()

# Insert a new, empty column before column N. This pushes the existing N, O
# and P columns (and all their data/styles) one column to the right, to O, P
# and Q respectively. The new column inherits the width of the column to its
# left (M), like Excel does when inserting a column.
$leftWidth = $ws.Columns("M").ColumnWidth
$ws.Range("N1").EntireColumn.Insert()
$ws.Columns("N").ColumnWidth = $leftWidth

# Update the selection on the sheet.
[void]$ws.Range("K13").Select()
